$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(53,8).Value = 178.61111  # H53: was 204.72974
$ws.Cells.Item(53,9).Value = 189.9375  # I53: was 202.4375
$ws.Cells.Item(53,10).Value = 169.55  # J53: was 206.4762
$ws.Cells.Item(53,11).Value = 189.9375  # K53: was 202.4375
$ws.Cells.Item(53,12).Value = 169.55  # L53: was 206.4762
$ws.Cells.Item(53,13).Value = 447.0625  # M53: was 434.5625
$ws.Cells.Item(53,14).Value = -1443.55  # N53: was -1480.4762
$ws.Cells.Item(70,8).Value = 44399.8  # H70: was 55250
$ws.Cells.Item(70,9).Value = 44399.8  # I70: was 55250
$ws.Cells.Item(70,11).Value = 133199.4  # K70: was 165750
$ws.Cells.Item(70,13).Value = -132929.4  # M70: was -165480
$ws.Cells.Item(73,8).Value = 44399.8  # H73: was 55250
$ws.Cells.Item(73,9).Value = 44399.8  # I73: was 55250
$ws.Cells.Item(73,11).Value = 133199.4  # K73: was 165750
$ws.Cells.Item(73,13).Value = -132263.4  # M73: was -164814
$ws.Cells.Item(74,8).Value = 5400  # H74: was 5500
$ws.Cells.Item(74,9).Value = 3000  # I74: was 4500
$ws.Cells.Item(74,10).Value = 6000  # J74: was 5666.6665
$ws.Cells.Item(74,11).Value = 3000  # K74: was 4500
$ws.Cells.Item(74,12).Value = 6000  # L74: was 5666.6665
$ws.Cells.Item(74,13).Value = -2064  # M74: was -3564
$ws.Cells.Item(74,14).Value = -7872  # N74: was -7538.6665
$ws.Cells.Item(77,8).Value = 5400  # H77: was 5500
$ws.Cells.Item(77,9).Value = 3000  # I77: was 4500
$ws.Cells.Item(77,10).Value = 6000  # J77: was 5666.6665
$ws.Cells.Item(77,11).Value = 15000  # K77: was 22500
$ws.Cells.Item(77,12).Value = 30000  # L77: was 28333.3325
$ws.Cells.Item(77,13).Value = -10320  # M77: was -17820
$ws.Cells.Item(77,14).Value = -39360  # N77: was -37693.3325
$ws.Cells.Item(97,8).Value = 1698.6666  # H97: was 1709.6666
$ws.Cells.Item(97,10).Value = 1698.6666  # J97: was 1709.6666
$ws.Cells.Item(97,12).Value = 5095.9998  # L97: was 5128.9998
$ws.Cells.Item(97,14).Value = -6087.9998  # N97: was -6120.9998
$ws.Cells.Item(98,8).Value = 2347.5715  # H98: was 2445.4285
$ws.Cells.Item(98,9).Value = 1766.6522  # I98: was 1885.7826
$ws.Cells.Item(98,11).Value = 1766.6522  # K98: was 1885.7826
$ws.Cells.Item(98,13).Value = -268.6522  # M98: was -387.7826
$ws.Cells.Item(116,8).Value = 16454.818  # H116: was 15500.25
$ws.Cells.Item(116,9).Value = 20001.25  # I116: was 17001
$ws.Cells.Item(116,11).Value = 20001.25  # K116: was 17001
$ws.Cells.Item(116,13).Value = -16559.25  # M116: was -13559
$ws.Cells.Item(121,8).Value = 1680.4736  # H121: was 1604.6666
$ws.Cells.Item(121,10).Value = 1680.4736  # J121: was 1604.6666
$ws.Cells.Item(121,12).Value = 5041.4208  # L121: was 4813.9998
$ws.Cells.Item(121,14).Value = -8535.4208  # N121: was -8307.9998
$ws.Cells.Item(122,8).Value = 2347.5715  # H122: was 2445.4285
$ws.Cells.Item(122,9).Value = 1766.6522  # I122: was 1885.7826
$ws.Cells.Item(122,11).Value = 5299.9566  # K122: was 5657.3478
$ws.Cells.Item(122,13).Value = -2849.9566  # M122: was -3207.3478
$ws.Cells.Item(125,8).Value = 10613.03  # H125: was 10930.9375
$ws.Cells.Item(125,10).Value = 1123.4286  # J125: was 1176
$ws.Cells.Item(125,12).Value = 10110.8574  # L125: was 10584
$ws.Cells.Item(125,14).Value = -15030.8574  # N125: was -15504
$ws.Cells.Item(137,8).Value = 7707.16  # H137: was 8268.817999999999
$ws.Cells.Item(137,9).Value = 9192.444  # I137: was 11823.833
$ws.Cells.Item(137,10).Value = 6871.6875  # J137: was 6935.6875
$ws.Cells.Item(137,11).Value = 27577.332  # K137: was 35471.499
$ws.Cells.Item(137,12).Value = 20615.0625  # L137: was 20807.0625
$ws.Cells.Item(137,13).Value = -25027.332  # M137: was -32921.499
$ws.Cells.Item(137,14).Value = -25715.0625  # N137: was -25907.0625
$ws.Cells.Item(138,8).Value = 2030.6346  # H138: was 2064.1296
$ws.Cells.Item(138,10).Value = 3504.2104  # J138: was 3450
$ws.Cells.Item(138,12).Value = 10512.6312  # L138: was 10350
$ws.Cells.Item(138,14).Value = -20792.6312  # N138: was -20630
$ws.Cells.Item(141,8).Value = 22137.93  # H141: was 23194.732
$ws.Cells.Item(141,9).Value = 23193.928  # I141: was 24359.076
$ws.Cells.Item(141,11).Value = 69581.784  # K141: was 73077.228
$ws.Cells.Item(141,13).Value = -64401.784  # M141: was -67897.228

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32,8).Value = 8042305  # H32: was 8155575
$ws.Cells.Item(32,9).Value = 1368480.9  # I32: was 1392915.8
$ws.Cells.Item(32,11).Value = 1368480.9  # K32: was 1392915.8
$ws.Cells.Item(32,13).Value = -1368193.9  # M32: was -1392628.8
$ws.Cells.Item(61,8).Value = 1887.4762  # H61: was 1921.85
$ws.Cells.Item(61,9).Value = 1826.8422  # I61: was 1861.6666
$ws.Cells.Item(61,11).Value = 1826.8422  # K61: was 1861.6666
$ws.Cells.Item(61,13).Value = -1614.8422  # M61: was -1649.6666
$ws.Cells.Item(122,8).Value = 7151.6665  # H122: was 7461.744
$ws.Cells.Item(122,9).Value = 8347.424000000001  # I122: was 8593.125
$ws.Cells.Item(122,10).Value = 3863.3333  # J122: was 4170.4546
$ws.Cells.Item(122,11).Value = 25042.272  # K122: was 25779.375
$ws.Cells.Item(122,12).Value = 11589.9999  # L122: was 12511.3638
$ws.Cells.Item(122,13).Value = -22592.272  # M122: was -23329.375
$ws.Cells.Item(122,14).Value = -16489.9999  # N122: was -17411.3638
$ws.Cells.Item(132,8).Value = 3070.9824  # H132: was 3223.3455
$ws.Cells.Item(132,9).Value = 2940.92  # I132: was 3110.0833
$ws.Cells.Item(132,11).Value = 8822.76  # K132: was 9330.249899999999
$ws.Cells.Item(132,13).Value = -6292.76  # M132: was -6800.249899999999
$ws.Cells.Item(136,8).Value = 1887.4762  # H136: was 1921.85
$ws.Cells.Item(136,9).Value = 1826.8422  # I136: was 1861.6666
$ws.Cells.Item(136,11).Value = 5480.5266  # K136: was 5584.9998
$ws.Cells.Item(136,13).Value = -2930.5266  # M136: was -3034.9998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(138,8).Value = 53166.668  # H138: was 0
$ws.Cells.Item(138,10).Value = 53166.668  # J138: was 0
$ws.Cells.Item(138,12).Value = 53166.668  # L138: was 0
$ws.Cells.Item(138,14).Value = -63446.668  # N138: was <<ABSENT>>

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22,8).Value = 163.81818  # H22: was 167.2
$ws.Cells.Item(22,9).Value = 170.75  # I22: was 176.57143
$ws.Cells.Item(22,11).Value = 170.75  # K22: was 176.57143
$ws.Cells.Item(22,13).Value = 179.25  # M22: was 173.42857
$ws.Cells.Item(31,8).Value = 5187.778  # H31: was 5337.3887
$ws.Cells.Item(31,9).Value = 2237  # I31: was 2685.8333
$ws.Cells.Item(31,11).Value = 2237  # K31: was 2685.8333
$ws.Cells.Item(31,13).Value = -1942  # M31: was -2390.8333
$ws.Cells.Item(34,8).Value = 5187.778  # H34: was 5337.3887
$ws.Cells.Item(34,9).Value = 2237  # I34: was 2685.8333
$ws.Cells.Item(34,11).Value = 2237  # K34: was 2685.8333
$ws.Cells.Item(34,13).Value = -2035  # M34: was -2483.8333
$ws.Cells.Item(58,8).Value = 1250  # H58: was 1006.61536
$ws.Cells.Item(58,9).Value = 1000  # I58: was 966
$ws.Cells.Item(58,10).Value = 1500  # J58: was 1230
$ws.Cells.Item(58,11).Value = 1000  # K58: was 966
$ws.Cells.Item(58,12).Value = 1500  # L58: was 1230
$ws.Cells.Item(58,13).Value = -797  # M58: was -763
$ws.Cells.Item(58,14).Value = -1906  # N58: was -1636
$ws.Cells.Item(88,8).Value = 24381  # H88: was 24621.5
$ws.Cells.Item(88,10).Value = 24381  # J88: was 24621.5
$ws.Cells.Item(88,12).Value = 24381  # L88: was 24621.5
$ws.Cells.Item(88,14).Value = -25193  # N88: was -25433.5
$ws.Cells.Item(91,8).Value = 24381  # H91: was 24621.5
$ws.Cells.Item(91,10).Value = 24381  # J91: was 24621.5
$ws.Cells.Item(91,12).Value = 24381  # L91: was 24621.5
$ws.Cells.Item(91,14).Value = -27189  # N91: was -27429.5
$ws.Cells.Item(122,8).Value = 3206  # H122: was 2803.6667
$ws.Cells.Item(122,9).Value = 3012  # I122: was 2505.5
$ws.Cells.Item(122,11).Value = 9036  # K122: was 7516.5
$ws.Cells.Item(122,13).Value = -6586  # M122: was -5066.5
$ws.Cells.Item(132,8).Value = 4628  # H132: was 4724.8335
$ws.Cells.Item(132,9).Value = 4686.1816  # I132: was 4828.222
$ws.Cells.Item(132,11).Value = 14058.5448  # K132: was 14484.666
$ws.Cells.Item(132,13).Value = -11528.5448  # M132: was -11954.666
$ws.Cells.Item(134,8).Value = 2210.0386  # H134: was 2332.6667
$ws.Cells.Item(134,9).Value = 1638.0476  # I134: was 1732.7368
$ws.Cells.Item(134,11).Value = 4914.142800000001  # K134: was 5198.2104
$ws.Cells.Item(134,13).Value = -2379.142800000001  # M134: was -2663.2104
$ws.Cells.Item(136,8).Value = 1250  # H136: was 1006.61536
$ws.Cells.Item(136,9).Value = 1000  # I136: was 966
$ws.Cells.Item(136,10).Value = 1500  # J136: was 1230
$ws.Cells.Item(136,11).Value = 3000  # K136: was 2898
$ws.Cells.Item(136,12).Value = 4500  # L136: was 3690
$ws.Cells.Item(136,13).Value = -450  # M136: was -348
$ws.Cells.Item(136,14).Value = -9600  # N136: was -8790
$ws.Cells.Item(141,8).Value = 89999.95  # H141: was 90000
$ws.Cells.Item(141,10).Value = 89999.95  # J141: was 90000
$ws.Cells.Item(141,12).Value = 89999.95  # L141: was 90000
$ws.Cells.Item(141,14).Value = -100359.95  # N141: was -100360

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5,8).Value = 1657.3889  # H5: was 1835.45
$ws.Cells.Item(5,9).Value = 615.8333  # I5: was 617.4545000000001
$ws.Cells.Item(5,10).Value = 3740.5  # J5: was 3324.111
$ws.Cells.Item(5,11).Value = 1847.4999  # K5: was 1852.3635
$ws.Cells.Item(5,12).Value = 11221.5  # L5: was 9972.332999999999
$ws.Cells.Item(5,13).Value = -1735.4999  # M5: was -1740.3635
$ws.Cells.Item(5,14).Value = -11445.5  # N5: was -10196.333
$ws.Cells.Item(68,8).Value = 2287.9048  # H68: was 2261.2727
$ws.Cells.Item(68,10).Value = 2258.2307  # J68: was 2218.5
$ws.Cells.Item(68,12).Value = 6774.6921  # L68: was 6655.5
$ws.Cells.Item(68,14).Value = -8396.6921  # N68: was -8277.5
$ws.Cells.Item(71,8).Value = 2287.9048  # H71: was 2261.2727
$ws.Cells.Item(71,10).Value = 2258.2307  # J71: was 2218.5
$ws.Cells.Item(71,12).Value = 20324.0763  # L71: was 19966.5
$ws.Cells.Item(71,14).Value = -28436.0763  # N71: was -28078.5
$ws.Cells.Item(107,8).Value = 993.875  # H107: was 976.90625
$ws.Cells.Item(107,9).Value = 527.6  # I107: was 498
$ws.Cells.Item(107,10).Value = 1080.2222  # J107: was 1087.4231
$ws.Cells.Item(107,11).Value = 1582.8  # K107: was 1494
$ws.Cells.Item(107,12).Value = 3240.6666  # L107: was 3262.2693
$ws.Cells.Item(107,13).Value = 337.1999999999998  # M107: was 426
$ws.Cells.Item(107,14).Value = -7080.6666  # N107: was -7102.2693
$ws.Cells.Item(135,8).Value = 1657.3889  # H135: was 1835.45
$ws.Cells.Item(135,9).Value = 615.8333  # I135: was 617.4545000000001
$ws.Cells.Item(135,10).Value = 3740.5  # J135: was 3324.111
$ws.Cells.Item(135,11).Value = 5542.4997  # K135: was 5557.0905
$ws.Cells.Item(135,12).Value = 33664.5  # L135: was 29916.999
$ws.Cells.Item(135,13).Value = -3007.4997  # M135: was -3022.0905
$ws.Cells.Item(135,14).Value = -38734.5  # N135: was -34986.999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2,8).Value = 565.3570999999999  # H2: was 624.5
$ws.Cells.Item(2,9).Value = 627.25  # I2: was 713.8570999999999
$ws.Cells.Item(2,10).Value = 482.83334  # J2: was 499.4
$ws.Cells.Item(2,11).Value = 627.25  # K2: was 713.8570999999999
$ws.Cells.Item(2,12).Value = 482.83334  # L2: was 499.4
$ws.Cells.Item(2,13).Value = -514.25  # M2: was -600.8570999999999
$ws.Cells.Item(2,14).Value = -708.83334  # N2: was -725.4
$ws.Cells.Item(102,8).Value = 2054.1292  # H102: was 2142.3447
$ws.Cells.Item(102,9).Value = 2081.0356  # I102: was 2181.5
$ws.Cells.Item(102,11).Value = 2081.0356  # K102: was 2181.5
$ws.Cells.Item(102,13).Value = -459.0356000000002  # M102: was -559.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(50,8).Value = 0  # H50: was 20000
$ws.Cells.Item(50,9).Value = 0  # I50: was 20000
$ws.Cells.Item(50,11).Value = 0  # K50: was 20000
$ws.Cells.Item(50,13).ClearContents()  # M50: was -19363
$ws.Cells.Item(100,8).Value = 26469.564  # H100: was 29154.125
$ws.Cells.Item(100,9).Value = 19024.586  # I100: was 20816.371
$ws.Cells.Item(100,11).Value = 19024.586  # K100: was 20816.371
$ws.Cells.Item(100,13).Value = -18483.586  # M100: was -20275.371
$ws.Cells.Item(122,8).Value = 5000  # H122: was 3333.3333
$ws.Cells.Item(122,9).Value = 0  # I122: was 1500
$ws.Cells.Item(122,10).Value = 5000  # J122: was 4250
$ws.Cells.Item(122,11).Value = 0  # K122: was 4500
$ws.Cells.Item(122,12).Value = 15000  # L122: was 12750
$ws.Cells.Item(122,13).ClearContents()  # M122: was -2050
$ws.Cells.Item(122,14).Value = -19900  # N122: was -17650
$ws.Cells.Item(132,8).Value = 2738.2334  # H132: was 2897.6072
$ws.Cells.Item(132,9).Value = 2611.8635  # I132: was 2778.3809
$ws.Cells.Item(132,10).Value = 3085.75  # J132: was 3255.2856
$ws.Cells.Item(132,11).Value = 7835.5905  # K132: was 8335.1427
$ws.Cells.Item(132,12).Value = 9257.25  # L132: was 9765.856800000001
$ws.Cells.Item(132,13).Value = -5305.5905  # M132: was -5805.1427
$ws.Cells.Item(132,14).Value = -14317.25  # N132: was -14825.8568
$ws.Cells.Item(136,8).Value = 2492.16  # H136: was 3067.9
$ws.Cells.Item(136,9).Value = 1539.6842  # I136: was 2013.4286
$ws.Cells.Item(136,10).Value = 5508.3335  # J136: was 5528.3335
$ws.Cells.Item(136,11).Value = 4619.0526  # K136: was 6040.2858
$ws.Cells.Item(136,12).Value = 16525.0005  # L136: was 16585.0005
$ws.Cells.Item(136,13).Value = -2069.0526  # M136: was -3490.2858
$ws.Cells.Item(136,14).Value = -21625.0005  # N136: was -21685.0005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81,8).Value = 2500750  # H81: was 1250874.8
$ws.Cells.Item(81,9).Value = 2500750  # I81: was 1250874.8
$ws.Cells.Item(81,11).Value = 5001500  # K81: was 2501749.6
$ws.Cells.Item(81,13).Value = -5000439  # M81: was -2500688.6
$ws.Cells.Item(84,8).Value = 2500750  # H84: was 1250874.8
$ws.Cells.Item(84,9).Value = 2500750  # I84: was 1250874.8
$ws.Cells.Item(84,11).Value = 25007500  # K84: was 12508748
$ws.Cells.Item(84,13).Value = -25002196  # M84: was -12503444
$ws.Cells.Item(132,8).Value = 11748.2  # H132: was 12539.143
$ws.Cells.Item(132,9).Value = 20796.125  # I132: was 23656
$ws.Cells.Item(132,10).Value = 1407.7142  # J132: was 1422.2858
$ws.Cells.Item(132,11).Value = 62388.375  # K132: was 70968
$ws.Cells.Item(132,12).Value = 4223.142599999999  # L132: was 4266.857400000001
$ws.Cells.Item(132,13).Value = -59858.375  # M132: was -68438
$ws.Cells.Item(132,14).Value = -9283.142599999999  # N132: was -9326.857400000001
